# Update cryptocurrency price/volume data per Feb 8 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.552.43"
$ws.Range("E2").Value = "  +3.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.431.40"
$ws.Range("E3").Value = "  +2.33%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.19"
$ws.Range("E5").Value = "  +3.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.00"
$ws.Range("E6").Value = "  +5.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.513"
$ws.Range("E7").Value = "  +1.72%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.509"
$ws.Range("E9").Value = "  +1.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.43"
$ws.Range("E10").Value = "  +3.60%  "

$ws.Range("E11").Value = "  +1.45%  "

$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("E13").Value = "  +2.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.98"
$ws.Range("E14").Value = "  +2.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.811.99"
$ws.Range("E15").Value = "  +2.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.425.79"
$ws.Range("E16").Value = "  +1.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.840"
$ws.Range("E17").Value = "  +4.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.472.73"
$ws.Range("E18").Value = "  +3.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.48"
$ws.Range("E19").Value = "  +2.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.42"
$ws.Range("E20").Value = "  +1.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("E21").Value = "  +2.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.92"
$ws.Range("E22").Value = "  +1.21%  "

$ws.Range("E23").Value = "  +3.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.79"
$ws.Range("E24").Value = "  +2.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.50"
$ws.Range("E25").Value = "  +2.68%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.22"
$ws.Range("E27").Value = "  +1.25%  "

$ws.Range("E28").Value = "  -2.94%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.69"
$ws.Range("E29").Value = "  +4.72%  "

$ws.Range("E30").Value = "  +5.58%  "

$ws.Range("E31").Value = "  +16.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.57"
$ws.Range("E32").Value = "  +12.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.20"
$ws.Range("E33").Value = "  +2.65%  "

$ws.Range("E34").Value = "  +0.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0766"
$ws.Range("E35").Value = "  +3.77%  "

$ws.Range("E36").Value = "  +2.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.51"
$ws.Range("E37").Value = "  +3.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.90"
$ws.Range("E38").Value = "  +3.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.33"
$ws.Range("E39").Value = "  +0.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "126.56"
$ws.Range("E40").Value = "  +11.88%  "

$ws.Range("E41").Value = "  +0.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.97"
$ws.Range("E42").Value = "  -3.77%  "

$ws.Range("E43").Value = "  +3.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.948.20"
$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.94"
$ws.Range("E46").Value = "  +7.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.58"
$ws.Range("E47").Value = "  +4.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.67"
$ws.Range("E48").Value = "  +10.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.57"
$ws.Range("E49").Value = "  +2.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.92"
$ws.Range("E50").Value = "  +2.18%  "

$ws.Range("E51").Value = "  +4.72%  "
